{"js": "// Replace each unique cell value with its updated value.\n// Each (old -> new) pair is unique in the source document, so a\n// simple exact search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"2025-08-15 Friday\", \"2025-08-16 Saturday\"],\n  [\"744\u00d76=\", \"137\u00d76=\"],\n  [\"712\u00d77=\", \"470\u00d74=\"],\n  [\"564\u00d78=\", \"362\u00d74=\"],\n  [\"824\u00d72=\", \"256\u00d73=\"],\n  [\"975\u00d78=\", \"472\u00d75=\"],\n  [\"789\u00d73=\", \"419\u00d78=\"],\n  [\"624\u00d76=\", \"940\u00d72=\"],\n  [\"535\u00d79=\", \"618\u00d77=\"],\n  [\"638\u00d76=\", \"758\u00d73=\"],\n  [\"449\u00d79=\", \"918\u00d73=\"],\n  [\"316\u00d79=\", \"127\u00d74=\"],\n  [\"190\u00d76=\", \"514\u00d76=\"],\n  [\"416\u00d72=\", \"259\u00d75=\"],\n  [\"925\u00d78=\", \"138\u00d78=\"],\n  [\"617\u00d74=\", \"512\u00d78=\"],\n  [\"482\u00d73=\", \"878\u00d78=\"],\n  [\"625\u00d76=\", \"272\u00d78=\"],\n  [\"462\u00d76=\", \"470\u00d75=\"],\n  [\"978\u00d74=\", \"127\u00d74=\"],\n  [\"797\u00d76=\", \"614\u00d79=\"],\n  [\"169\u00d79=\", \"848\u00d78=\"],\n  [\"500\u00d74=\", \"185\u00d76=\"],\n  [\"436\u00d76=\", \"606\u00d78=\"],\n  [\"206\u00d74=\", \"755\u00d78=\"],\n  [\"530\u00d76=\", \"927\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for '\" + oldText + \"' but found \" + results.items.length\n    );\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n", "ps1": "# Replace each unique cell value with its updated value using\n# Word's Find/Replace. Each (old -> new) pair is unique in the\n# source document, so an exact, non-wildcard Find.Execute with\n# Replace:=wdReplaceAll (2) is unambiguous for every pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-08-15 Friday\", \"2025-08-16 Saturday\"),\n    @(\"744\u00d76=\", \"137\u00d76=\"),\n    @(\"712\u00d77=\", \"470\u00d74=\"),\n    @(\"564\u00d78=\", \"362\u00d74=\"),\n    @(\"824\u00d72=\", \"256\u00d73=\"),\n    @(\"975\u00d78=\", \"472\u00d75=\"),\n    @(\"789\u00d73=\", \"419\u00d78=\"),\n    @(\"624\u00d76=\", \"940\u00d72=\"),\n    @(\"535\u00d79=\", \"618\u00d77=\"),\n    @(\"638\u00d76=\", \"758\u00d73=\"),\n    @(\"449\u00d79=\", \"918\u00d73=\"),\n    @(\"316\u00d79=\", \"127\u00d74=\"),\n    @(\"190\u00d76=\", \"514\u00d76=\"),\n    @(\"416\u00d72=\", \"259\u00d75=\"),\n    @(\"925\u00d78=\", \"138\u00d78=\"),\n    @(\"617\u00d74=\", \"512\u00d78=\"),\n    @(\"482\u00d73=\", \"878\u00d78=\"),\n    @(\"625\u00d76=\", \"272\u00d78=\"),\n    @(\"462\u00d76=\", \"470\u00d75=\"),\n    @(\"978\u00d74=\", \"127\u00d74=\"),\n    @(\"797\u00d76=\", \"614\u00d79=\"),\n    @(\"169\u00d79=\", \"848\u00d78=\"),\n    @(\"500\u00d74=\", \"185\u00d76=\"),\n    @(\"436\u00d76=\", \"606\u00d78=\"),\n    @(\"206\u00d74=\", \"755\u00d78=\"),\n    @(\"530\u00d76=\", \"927\u00d76=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue: keep searching the whole story\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $found = $find.Execute(\n        $oldText,     # FindText\n        $true,        # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap = wdFindContinue\n        $false,       # Format\n        $newText,     # ReplaceWith\n        2             # Replace = wdReplaceAll\n    )\n\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n\n"}
